# Apply crypto price/volume updates per the commit diff.
# Numeric-looking Price/Volume strings get a leading "'" (apostrophe)
# so Excel stores them as text verbatim instead of coercing to numbers
# (preserving things like trailing zeros and dotted thousand separators).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.041.60"
$ws.Range("E2").Value = "'  -0.07%  "

$ws.Range("D3").Value = "'2.302.71"
$ws.Range("E3").Value = "'  +0.03%  "

$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'300.29"
$ws.Range("E5").Value = "'  -0.24%  "

$ws.Range("D6").Value = "'97.98"
$ws.Range("E6").Value = "'  -1.66%  "

$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "'  +3.10%  "

$ws.Range("E8").Value = "'  -0.02%  "

$ws.Range("E9").Value = "'  +1.24%  "

$ws.Range("D10").Value = "'36.28"
$ws.Range("E10").Value = "'  -0.17%  "

$ws.Range("E11").Value = "'  +0.24%  "

$ws.Range("E12").Value = "'  +0.60%  "

$ws.Range("D13").Value = "'17.72"
$ws.Range("E13").Value = "'  -2.50%  "

$ws.Range("E14").Value = "'  -0.60%  "

$ws.Range("D15").Value = "'2.661.18"
$ws.Range("E15").Value = "'  +0.15%  "

$ws.Range("D16").Value = "'2.336.03"
$ws.Range("E16").Value = "'  +3.27%  "

$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = "'  -1.20%  "

$ws.Range("D18").Value = "'42.933.49"
$ws.Range("E18").Value = "'  -0.10%  "

$ws.Range("D19").Value = "'12.96"
$ws.Range("E19").Value = "'  +2.93%  "

$ws.Range("E20").Value = "'  +0.76%  "

$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "'  +0.23%  "

$ws.Range("D22").Value = "'68.30"
$ws.Range("E22").Value = "'  +0.55%  "

$ws.Range("D23").Value = "'238.10"
$ws.Range("E23").Value = "'  +1.01%  "

$ws.Range("E24").Value = "'  -1.43%  "

$ws.Range("E25").Value = "'  -0.72%  "

$ws.Range("E26").Value = "'  -0.50%  "

$ws.Range("D27").Value = "'4.01"
$ws.Range("E27").Value = "'  -0.47%  "

$ws.Range("D28").Value = "'24.99"
$ws.Range("E28").Value = "'  +0.00%  "

$ws.Range("D29").Value = "'2.05"
$ws.Range("E29").Value = "'  -12.90%  "

$ws.Range("D30").Value = "'9.17"
$ws.Range("E30").Value = "'  +0.28%  "

$ws.Range("D31").Value = "'163.42"

$ws.Range("D32").Value = "'33.17"
$ws.Range("E32").Value = "'  -4.28%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  -0.01%  "

$ws.Range("D34").Value = "'5.13"
$ws.Range("E34").Value = "'  +2.09%  "

$ws.Range("D35").Value = "'18.17"
$ws.Range("E35").Value = "'  +2.58%  "

$ws.Range("D36").Value = "'4.80"
$ws.Range("E36").Value = "'  +3.97%  "

$ws.Range("E37").Value = "'  +0.34%  "

$ws.Range("D38").Value = "'0.0698"
$ws.Range("E38").Value = "'  +1.15%  "

$ws.Range("E39").Value = "'  +0.60%  "

$ws.Range("E40").Value = "'  -0.54%  "

$ws.Range("D41").Value = "'2.79"
$ws.Range("E41").Value = "'  -1.26%  "

$ws.Range("E42").Value = "'  +1.32%  "

$ws.Range("D43").Value = "'2.016.46"
$ws.Range("E43").Value = "'  +1.89%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.27"
$ws.Range("E44").Value = "'  -1.76%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0287"
$ws.Range("E45").Value = "'  -1.48%  "

$ws.Range("D46").Value = "'10.35"
$ws.Range("E46").Value = "'  +1.59%  "

$ws.Range("D47").Value = "'17.47"
$ws.Range("E47").Value = "'  -0.09%  "

$ws.Range("E48").Value = "'  -2.34%  "

$ws.Range("D49").Value = "'54.38"
$ws.Range("E49").Value = "'  -2.28%  "

$ws.Range("D50").Value = "'2.526.57"
$ws.Range("E50").Value = "'  +0.17%  "

$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "'  -1.02%  "
